# Update "想去人数" (attendance count) figures for the 丽水-漫展信息 workbook.
# Both the "展览" and "全部类型" sheets carry the same table and both need
# the same three cells updated:
#   F2: 104 -> 107
#   F3: 940 -> 941
#   F5: 429 -> 430

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 107
    $ws.Range("F3").Value = 941
    $ws.Range("F5").Value = 430
}
